$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BRENNER's balance (row 2, column C)
$ws.Cells.Item(2, 3).Value = 63764.66

# Delete rows for ADELE (row 3), MARCUS (row 6), NATALIA (row 7)
# Delete from bottom to top to keep row indices valid
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(3).Delete()
